$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 69: "H69" = 9500, "J69" = 9500, "L69" = 28500, "N69" = -30248
$ws.Range("H69").Value = 9500
$ws.Range("J69").Value = 9500
$ws.Range("L69").Value = 28500
$ws.Range("N69").Value = -30248
# Row 72: "H72" = 9500, "J72" = 9500, "L72" = 85500, "N72" = -94236
$ws.Range("H72").Value = 9500
$ws.Range("J72").Value = 9500
$ws.Range("L72").Value = 85500
$ws.Range("N72").Value = -94236
# Row 137: "H137" = 6505.4736, "I137" = 2275.3, "J137" = 11205.667, "K137" = 6825.900000000001, "L137" = 33617.001, "M137" = -4275.900000000001, "N137" = -38717.001
$ws.Range("H137").Value = 6505.4736
$ws.Range("I137").Value = 2275.3
$ws.Range("J137").Value = 11205.667
$ws.Range("K137").Value = 6825.900000000001
$ws.Range("L137").Value = 33617.001
$ws.Range("M137").Value = -4275.900000000001
$ws.Range("N137").Value = -38717.001
# Row 138: "H138" = 5469.625, "J138" = 6250.346, "L138" = 18751.038, "N138" = -29031.038
$ws.Range("H138").Value = 5469.625
$ws.Range("J138").Value = 6250.346
$ws.Range("L138").Value = 18751.038
$ws.Range("N138").Value = -29031.038

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2: "H2" = 1533.4445, "I2" = 1621.0667, "K2" = 1621.0667, "M2" = -1508.0667
$ws.Range("H2").Value = 1533.4445
$ws.Range("I2").Value = 1621.0667
$ws.Range("K2").Value = 1621.0667
$ws.Range("M2").Value = -1508.0667
# Row 32: "H32" = 16312.028, "I32" = 6767.425, "J32" = 28627.645, "K32" = 6767.425, "L32" = 28627.645, "M32" = -6480.425, "N32" = -29201.645
$ws.Range("H32").Value = 16312.028
$ws.Range("I32").Value = 6767.425
$ws.Range("J32").Value = 28627.645
$ws.Range("K32").Value = 6767.425
$ws.Range("L32").Value = 28627.645
$ws.Range("M32").Value = -6480.425
$ws.Range("N32").Value = -29201.645
# Row 45: "H45" = 2176.4546, "I45" = 1490.1666, "K45" = 1490.1666, "M45" = -1113.1666
$ws.Range("H45").Value = 2176.4546
$ws.Range("I45").Value = 1490.1666
$ws.Range("K45").Value = 1490.1666
$ws.Range("M45").Value = -1113.1666
# Row 61: "H61" = 2418.3, "I61" = 2418.3, "K61" = 2418.3, "M61" = -2206.3
$ws.Range("H61").Value = 2418.3
$ws.Range("I61").Value = 2418.3
$ws.Range("K61").Value = 2418.3
$ws.Range("M61").Value = -2206.3
# Row 74: "H74" = 3815.182, "I74" = 1501.8462, "J74" = 7156.6665, "K74" = 1501.8462, "L74" = 7156.6665, "M74" = -627.8462, "N74" = -8904.666499999999
$ws.Range("H74").Value = 3815.182
$ws.Range("I74").Value = 1501.8462
$ws.Range("J74").Value = 7156.6665
$ws.Range("K74").Value = 1501.8462
$ws.Range("L74").Value = 7156.6665
$ws.Range("M74").Value = -627.8462
$ws.Range("N74").Value = -8904.666499999999
# Row 77: "H77" = 3815.182, "I77" = 1501.8462, "J77" = 7156.6665, "K77" = 7509.231, "L77" = 35783.3325, "M77" = -3141.231, "N77" = -44519.3325
$ws.Range("H77").Value = 3815.182
$ws.Range("I77").Value = 1501.8462
$ws.Range("J77").Value = 7156.6665
$ws.Range("K77").Value = 7509.231
$ws.Range("L77").Value = 35783.3325
$ws.Range("M77").Value = -3141.231
$ws.Range("N77").Value = -44519.3325
# Row 80: "H80" = 89999, "J80" = 89999, "L80" = 89999, "N80" = -91995
$ws.Range("H80").Value = 89999
$ws.Range("J80").Value = 89999
$ws.Range("L80").Value = 89999
$ws.Range("N80").Value = -91995
# Row 83: "H83" = 89999, "J83" = 89999, "L83" = 269997, "N83" = -279981
$ws.Range("H83").Value = 89999
$ws.Range("J83").Value = 89999
$ws.Range("L83").Value = 269997
$ws.Range("N83").Value = -279981
# Row 97: "H97" = 793.125, "I97" = 863, "J97" = 676.6667, "K97" = 863, "L97" = 676.6667, "M97" = -367, "N97" = -1668.6667
$ws.Range("H97").Value = 793.125
$ws.Range("I97").Value = 863
$ws.Range("J97").Value = 676.6667
$ws.Range("K97").Value = 863
$ws.Range("L97").Value = 676.6667
$ws.Range("M97").Value = -367
$ws.Range("N97").Value = -1668.6667
# Row 116: "H116" = 1533.4445, "I116" = 1621.0667, "K116" = 1621.0667, "M116" = 672.9332999999999
$ws.Range("H116").Value = 1533.4445
$ws.Range("I116").Value = 1621.0667
$ws.Range("K116").Value = 1621.0667
$ws.Range("M116").Value = 672.9332999999999
# Row 122: "H122" = 360352.53, "I122" = 558049.7, "K122" = 1674149.1, "M122" = -1671699.1
$ws.Range("H122").Value = 360352.53
$ws.Range("I122").Value = 558049.7
$ws.Range("K122").Value = 1674149.1
$ws.Range("M122").Value = -1671699.1
# Row 136: "H136" = 2418.3, "I136" = 2418.3, "K136" = 7254.900000000001, "M136" = -4704.900000000001
$ws.Range("H136").Value = 2418.3
$ws.Range("I136").Value = 2418.3
$ws.Range("K136").Value = 7254.900000000001
$ws.Range("M136").Value = -4704.900000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3: "H3" = 1533.4445, "I3" = 1621.0667, "K3" = 1621.0667, "M3" = -1507.0667
$ws.Range("H3").Value = 1533.4445
$ws.Range("I3").Value = 1621.0667
$ws.Range("K3").Value = 1621.0667
$ws.Range("M3").Value = -1507.0667
# Row 80: "H80" = 387, "I80" = 328.875, "K80" = 328.875, "M80" = 669.125
$ws.Range("H80").Value = 387
$ws.Range("I80").Value = 328.875
$ws.Range("K80").Value = 328.875
$ws.Range("M80").Value = 669.125
# Row 83: "H83" = 387, "I83" = 328.875, "K83" = 1644.375, "M83" = 3347.625
$ws.Range("H83").Value = 387
$ws.Range("I83").Value = 328.875
$ws.Range("K83").Value = 1644.375
$ws.Range("M83").Value = 3347.625
# Row 86: "H86" = 3378.1667, "J86" = 4403.5, "L86" = 4403.5, "N86" = -6649.5
$ws.Range("H86").Value = 3378.1667
$ws.Range("J86").Value = 4403.5
$ws.Range("L86").Value = 4403.5
$ws.Range("N86").Value = -6649.5
# Row 89: "H89" = 3378.1667, "J89" = 4403.5, "L89" = 22017.5, "N89" = -33249.5
$ws.Range("H89").Value = 3378.1667
$ws.Range("J89").Value = 4403.5
$ws.Range("L89").Value = 22017.5
$ws.Range("N89").Value = -33249.5
# Row 94: "H94" = 667.7826, "I94" = 636.85, "K94" = 636.85, "M94" = -185.85
$ws.Range("H94").Value = 667.7826
$ws.Range("I94").Value = 636.85
$ws.Range("K94").Value = 636.85
$ws.Range("M94").Value = -185.85
# Row 99: "H99" = 2088.8333
$ws.Range("H99").Value = 2088.8333

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31: "H31" = 4893.2173, "J31" = 5345.5386, "L31" = 5345.5386, "N31" = -5935.5386
$ws.Range("H31").Value = 4893.2173
$ws.Range("J31").Value = 5345.5386
$ws.Range("L31").Value = 5345.5386
$ws.Range("N31").Value = -5935.5386
# Row 34: "H34" = 4893.2173, "J34" = 5345.5386, "L34" = 5345.5386, "N34" = -5749.5386
$ws.Range("H34").Value = 4893.2173
$ws.Range("J34").Value = 5345.5386
$ws.Range("L34").Value = 5345.5386
$ws.Range("N34").Value = -5749.5386
# Row 58: "H58" = 4051.1304, "J58" = 7643.875, "L58" = 7643.875, "N58" = -8049.875
$ws.Range("H58").Value = 4051.1304
$ws.Range("J58").Value = 7643.875
$ws.Range("L58").Value = 7643.875
$ws.Range("N58").Value = -8049.875
# Row 99: "H99" = 17784.072, "I99" = 15398.4, "K99" = 15398.4, "M99" = -13900.4
$ws.Range("H99").Value = 17784.072
$ws.Range("I99").Value = 15398.4
$ws.Range("K99").Value = 15398.4
$ws.Range("M99").Value = -13900.4
# Row 126: "H126" = 17784.072, "I126" = 15398.4, "K126" = 46195.2, "M126" = -43725.2
$ws.Range("H126").Value = 17784.072
$ws.Range("I126").Value = 15398.4
$ws.Range("K126").Value = 46195.2
$ws.Range("M126").Value = -43725.2
# Row 132: "H132" = 1979.3549, "I132" = 1867.8148, "J132" = 2732.25, "K132" = 5603.4444, "L132" = 8196.75, "M132" = -3073.4444, "N132" = -13256.75
$ws.Range("H132").Value = 1979.3549
$ws.Range("I132").Value = 1867.8148
$ws.Range("J132").Value = 2732.25
$ws.Range("K132").Value = 5603.4444
$ws.Range("L132").Value = 8196.75
$ws.Range("M132").Value = -3073.4444
$ws.Range("N132").Value = -13256.75
# Row 136: "H136" = 4051.1304, "J136" = 7643.875, "L136" = 22931.625, "N136" = -28031.625
$ws.Range("H136").Value = 4051.1304
$ws.Range("J136").Value = 7643.875
$ws.Range("L136").Value = 22931.625
$ws.Range("N136").Value = -28031.625

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 34: "H34" = 2173.5293, "I34" = 1691.6666, "J34" = 2436.3635, "K34" = 5074.9998, "L34" = 7309.0905, "M34" = -4990.9998, "N34" = -7477.0905
$ws.Range("H34").Value = 2173.5293
$ws.Range("I34").Value = 1691.6666
$ws.Range("J34").Value = 2436.3635
$ws.Range("K34").Value = 5074.9998
$ws.Range("L34").Value = 7309.0905
$ws.Range("M34").Value = -4990.9998
$ws.Range("N34").Value = -7477.0905
# Row 131: "H131" = 1311.5, "J131" = 1465.5, "L131" = 4396.5, "N131" = -14476.5
$ws.Range("H131").Value = 1311.5
$ws.Range("J131").Value = 1465.5
$ws.Range("L131").Value = 4396.5
$ws.Range("N131").Value = -14476.5
# Row 134: "H134" = 2202.875, "I134" = 660.5714, "K134" = 1981.7142, "M134" = 3088.2858
$ws.Range("H134").Value = 2202.875
$ws.Range("I134").Value = 660.5714
$ws.Range("K134").Value = 1981.7142
$ws.Range("M134").Value = 3088.2858
# Row 137: "H137" = 3460.6667, "I137" = 2250, "J137" = 3702.8, "K137" = 6750, "L137" = 11108.4, "M137" = -1650, "N137" = -21308.4
$ws.Range("H137").Value = 3460.6667
$ws.Range("I137").Value = 2250
$ws.Range("J137").Value = 3702.8
$ws.Range("K137").Value = 6750
$ws.Range("L137").Value = 11108.4
$ws.Range("M137").Value = -1650
$ws.Range("N137").Value = -21308.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22: "H22" = 1661.125, "I22" = 998.3333, "J22" = 2058.8, "K22" = 998.3333, "L22" = 2058.8, "M22" = -703.3333, "N22" = -2648.8
$ws.Range("H22").Value = 1661.125
$ws.Range("I22").Value = 998.3333
$ws.Range("J22").Value = 2058.8
$ws.Range("K22").Value = 998.3333
$ws.Range("L22").Value = 2058.8
$ws.Range("M22").Value = -703.3333
$ws.Range("N22").Value = -2648.8
# Row 27: "H27" = 1661.125, "I27" = 998.3333, "J27" = 2058.8, "K27" = 998.3333, "L27" = 2058.8, "M27" = -891.3333, "N27" = -2272.8
$ws.Range("H27").Value = 1661.125
$ws.Range("I27").Value = 998.3333
$ws.Range("J27").Value = 2058.8
$ws.Range("K27").Value = 998.3333
$ws.Range("L27").Value = 2058.8
$ws.Range("M27").Value = -891.3333
$ws.Range("N27").Value = -2272.8
# Row 40: "H40" = 5500, "I40" = 5000, "J40" = 6000, "K40" = 5000, "L40" = 6000, "M40" = -4864, "N40" = -6272
$ws.Range("H40").Value = 5500
$ws.Range("I40").Value = 5000
$ws.Range("J40").Value = 6000
$ws.Range("K40").Value = 5000
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = -4864
$ws.Range("N40").Value = -6272
# Row 46: "H46" = 4038.3076, "I46" = 2720, "K46" = 2720, "M46" = -2532
$ws.Range("H46").Value = 4038.3076
$ws.Range("I46").Value = 2720
$ws.Range("K46").Value = 2720
$ws.Range("M46").Value = -2532
# Row 61: "H61" = 3769.963, "I61" = 3469.3044, "K61" = 3469.3044, "M61" = -3267.3044
$ws.Range("H61").Value = 3769.963
$ws.Range("I61").Value = 3469.3044
$ws.Range("K61").Value = 3469.3044
$ws.Range("M61").Value = -3267.3044
# Row 93: "H93" = 1284.2858, "I93" = 497.5, "K93" = 497.5, "M93" = 750.5
$ws.Range("H93").Value = 1284.2858
$ws.Range("I93").Value = 497.5
$ws.Range("K93").Value = 497.5
$ws.Range("M93").Value = 750.5
# Row 113: "H113" = 3769.963, "I113" = 3469.3044, "K113" = 3469.3044, "M113" = -1299.3044
$ws.Range("H113").Value = 3769.963
$ws.Range("I113").Value = 3469.3044
$ws.Range("K113").Value = 3469.3044
$ws.Range("M113").Value = -1299.3044
# Row 122: "H122" = 7614.8184, "I122" = 6073.25, "K122" = 18219.75, "M122" = -15769.75
$ws.Range("H122").Value = 7614.8184
$ws.Range("I122").Value = 6073.25
$ws.Range("K122").Value = 18219.75
$ws.Range("M122").Value = -15769.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122: "H122" = 2399.5, "I122" = 0, "J122" = 2399.5, "K122" = 0, "L122" = 7198.5, "N122" = -12098.5
$ws.Range("H122").Value = 2399.5
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2399.5
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 7198.5
$ws.Range("N122").Value = -12098.5
$ws.Range("M122").ClearContents()  # cell removed in target
# Row 132: "H132" = 2865.8948, "I132" = 2269.818, "J132" = 3685.5, "K132" = 6809.454000000001, "L132" = 11056.5, "M132" = -4279.454000000001, "N132" = -16116.5
$ws.Range("H132").Value = 2865.8948
$ws.Range("I132").Value = 2269.818
$ws.Range("J132").Value = 3685.5
$ws.Range("K132").Value = 6809.454000000001
$ws.Range("L132").Value = 11056.5
$ws.Range("M132").Value = -4279.454000000001
$ws.Range("N132").Value = -16116.5

Write-Output "Applied Seraph_Profits market-data refresh edits."